$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D3" = -6.964699999999993
    "C7" = -13.64769999999999
    "B8" = 6.076899999999999
    "B10" = 4.826399999999999
    "E10" = 15.9611
    "B12" = 5.106899999999999
    "E12" = 17.66640000000001
    "E13" = 16.3844
    "E14" = 16.6628
    "C15" = -13.99339999999998
    "B18" = 6.720899999999991
    "C18" = -13.70679999999999
    "D18" = -9.014900000000001
    "D19" = -9.128099999999991
    "C20" = -12.1341
    "D27" = -8.722600000000007
    "C29" = -11.7873
    "E29" = 17.27660000000002
    "C30" = -13.18989999999999
    "C31" = -12.71149999999999
    "D31" = -9.070500000000008
    "E32" = 15.966
    "E35" = 16.55619999999999
    "B37" = 8.557700000000006
    "D38" = -8.653100000000002
    "C40" = -13.2374
    "D42" = -9.011199999999995
    "E43" = 17.42610000000001
    "D44" = -7.651599999999999
    "D47" = -7.4461
    "E48" = 17.46390000000002
    "E49" = 15.52809999999999
    "C50" = -13.67949999999999
    "E50" = 16.5171
    "B55" = 6.140999999999994
    "E56" = 15.9898
    "D58" = -8.097399999999995
    "D65" = -7.6391
    "B68" = 5.488699999999999
    "C68" = -11.3648
    "E69" = 17.44810000000004
    "D73" = -7.721
    "C76" = -12.8677
    "B77" = 9.344900000000001
    "B78" = 9.365900000000003
    "B81" = 5.122000000000006
    "E81" = 16.0194
    "B82" = 6.010600000000002
    "C87" = -14.27039999999999
    "C88" = -13.79049999999999
    "D90" = -7.876900000000004
    "E92" = 18.52000000000002
    "D94" = -6.787699999999997
    "D95" = -7.765799999999993
    "C96" = -13.3234
    "C98" = -11.95999999999999
    "C101" = -13.8705
    "D101" = -8.043099999999997
    "C102" = -13.26880000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"